# "server notify new player"
# Append a new "to-do" row to the list: a new task describing a server
# notification when a player joins, marked as already done ("Cделано").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 19

# Copy formatting from existing rows that already use the same visual
# style so the new row blends in with the rest of the table:
#   - column A uses the wrapped/bordered "task text" style (row 1, col A)
#   - column B uses the green "Cделано" (done) style (row 1, col B)
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4122) # xlPasteFormats

# Fill in the new task text and its status.
$ws.Cells.Item($newRow, 1).Value = "Добавить оповещение о присоединении игрока"
$ws.Cells.Item($newRow, 2).Value = "Cделано"

# Match the row height used by the other two-line task rows.
$ws.Rows.Item($newRow).RowHeight = 30

# Move the visible selection the way the author left it after editing.
$ws.Range("A13").Select()
$win = $excel.ActiveWindow
if ($win -ne $null) {
    $win.ScrollRow = 13
    $win.ScrollColumn = 1
}
$ws.Range("B17").Select()
